# "Gantt à jour MATIN 27/10"
# Advance the Gantt chart's scroll position and refresh a handful of
# milestone figures in the "Jalons" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gantt")
$ws.Activate()

# --- Scroll the timeline forward (scrollbar linked to $E$3) -----------------
# Moving the "Incrément de défilement" scrollbar to 25 shifts the whole
# Gantt timeline (H5:BK5 and everything derived from it) forward by 25 days.
$ws.Range("E3").Value = 25

# Keep the scrollbar control itself (and the window's visible viewport) in
# sync with the new increment / selection, best effort.
$scrollBar = $ws.Shapes.Item(1)
$scrollBar.ControlFormat.Value = 25

$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1

# --- Update milestone progress / scheduling data -----------------------------
# Row 18 — "Rédaction de la documentation Développeur" (AMINE): progress up
# to 80%, start date slipped to 13/10, duration shortened to 16 days.
$ws.Range("D18").Value = 0.8
$ws.Range("E18").Value = 44117
$ws.Range("F18").Value = 16

# Row 23 — now spans 1 day instead of 0.
$ws.Range("F23").Value = 1

# Row 29 — "Portage du site internet ..." progress up to 80%.
$ws.Range("D29").Value = 0.8

# Row 30 — "Test d'utilisation" progress up to 60%.
$ws.Range("D30").Value = 0.6

# --- Move the active selection to reflect where work left off ---------------
$ws.Range("F24").Select() | Out-Null
